# C1--C2-and-C3-PowerPoint.pptx -- apply the commit's change.
#
# Slide 16 contains a 2-column "Cash flow" summary table whose built-in
# table style was changed (Table Styles gallery) from GUID
# {C57A9A4E-9213-427A-8EDB-795F2338FB74} ("Table_0") to
# {B0F290A6-BF9D-4C0D-AA03-E379C7105F1D}. Table styles cannot be
# assigned through the Table.Style property directly -- PowerPoint
# raises "Table styles cannot be assigned through a property - call
# Table.ApplyStyle("{GUID}") instead" -- so use Table.ApplyStyle with
# the target style's GUID, exactly like choosing a new style from the
# gallery would.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B0F290A6-BF9D-4C0D-AA03-E379C7105F1D}")
    }
}
